$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2465.1924
$ws.Range("I100").Value = 2411.4707
$ws.Range("J100").Value = 2566.6667
$ws.Range("K100").Value = 2411.4707
$ws.Range("L100").Value = 2566.6667
$ws.Range("M100").Value = -1870.4707
$ws.Range("N100").Value = -3648.6667

$ws.Range("H113").Value = 2352.7778
$ws.Range("I113").Value = 2156
$ws.Range("K113").Value = 2156
$ws.Range("M113").Value = 1098

$ws.Range("H121").Value = 3578.3333
$ws.Range("J121").Value = 3838.125
$ws.Range("L121").Value = 11514.375
$ws.Range("N121").Value = -15008.375

$ws.Range("H131").Value = 2131.8235
$ws.Range("J131").Value = 3369.6667
$ws.Range("L131").Value = 10109.0001
$ws.Range("N131").Value = -20189.0001

$ws.Range("H138").Value = 128859.91
$ws.Range("I138").Value = 1785.9259
$ws.Range("J138").Value = 187012.4
$ws.Range("K138").Value = 5357.7777
$ws.Range("L138").Value = 561037.2
$ws.Range("M138").Value = -217.7776999999996
$ws.Range("N138").Value = -571317.2

$ws.Range("H141").Value = 3539.5
$ws.Range("I141").Value = 1725.9259
$ws.Range("J141").Value = 19861.666
$ws.Range("K141").Value = 5177.7777
$ws.Range("L141").Value = 59584.99800000001
$ws.Range("M141").Value = 2.222300000000359
$ws.Range("N141").Value = -69944.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 40
$ws.Range("J14").Value = 40
$ws.Range("L14").Value = 40
$ws.Range("N14").Value = -390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 31252000
$ws.Range("I105").Value = 31252000
$ws.Range("K105").Value = 31252000
$ws.Range("M105").Value = -31250253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 379
$ws.Range("I22").Value = 273
$ws.Range("K22").Value = 273
$ws.Range("M22").Value = 77

$ws.Range("H31").Value = 6220.1763
$ws.Range("I31").Value = 3000.6
$ws.Range("J31").Value = 6775.276
$ws.Range("K31").Value = 3000.6
$ws.Range("L31").Value = 6775.276
$ws.Range("M31").Value = -2705.6
$ws.Range("N31").Value = -7365.276

$ws.Range("H34").Value = 6220.1763
$ws.Range("I34").Value = 3000.6
$ws.Range("J34").Value = 6775.276
$ws.Range("K34").Value = 3000.6
$ws.Range("L34").Value = 6775.276
$ws.Range("M34").Value = -2798.6
$ws.Range("N34").Value = -7179.276

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

$ws.Range("H69").Value = 25360
$ws.Range("I69").Value = 7400
$ws.Range("J69").Value = 37333.332
$ws.Range("K69").Value = 7400
$ws.Range("L69").Value = 37333.332
$ws.Range("M69").Value = -6651
$ws.Range("N69").Value = -38831.332

$ws.Range("H72").Value = 25360
$ws.Range("I72").Value = 7400
$ws.Range("J72").Value = 37333.332
$ws.Range("K72").Value = 22200
$ws.Range("L72").Value = 111999.996
$ws.Range("M72").Value = -18456
$ws.Range("N72").Value = -119487.996

$ws.Range("H105").Value = 535
$ws.Range("I105").Value = 502
$ws.Range("J105").Value = 700
$ws.Range("K105").Value = 502
$ws.Range("L105").Value = 700
$ws.Range("M105").Value = 1245
$ws.Range("N105").Value = -4194

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1454.1
$ws.Range("I18").Value = 2323.5
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = 6970.5
$ws.Range("L18").Value = 450
$ws.Range("M18").Value = -6801.5
$ws.Range("N18").Value = -788

$ws.Range("H44").Value = 500
$ws.Range("I44").Value = 500
$ws.Range("K44").Value = 1500
$ws.Range("M44").Value = -1102

$ws.Range("H68").Value = 12788
$ws.Range("I68").Value = 19320.8
$ws.Range("J68").Value = 1900
$ws.Range("K68").Value = 57962.39999999999
$ws.Range("L68").Value = 5700
$ws.Range("M68").Value = -57151.39999999999
$ws.Range("N68").Value = -7322

$ws.Range("H71").Value = 12788
$ws.Range("I71").Value = 19320.8
$ws.Range("J71").Value = 1900
$ws.Range("K71").Value = 173887.2
$ws.Range("L71").Value = 17100
$ws.Range("M71").Value = -169831.2
$ws.Range("N71").Value = -25212

$ws.Range("H97").Value = 14971.857
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 25750.75
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 77252.25
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -78244.25

$ws.Range("H104").Value = 5544.5
$ws.Range("J104").Value = 5544.5
$ws.Range("L104").Value = 16633.5
$ws.Range("N104").Value = -21875.5

$ws.Range("H106").Value = 8197.959999999999
$ws.Range("J106").Value = 8197.959999999999
$ws.Range("L106").Value = 24593.88
$ws.Range("N106").Value = -26485.88

$ws.Range("H113").Value = 937.1458
$ws.Range("I113").Value = 715.1923
$ws.Range("J113").Value = 1199.4546
$ws.Range("K113").Value = 2145.5769
$ws.Range("L113").Value = 3598.3638
$ws.Range("M113").Value = 24.42309999999998
$ws.Range("N113").Value = -7938.3638

$ws.Range("H122").Value = 7736.7144
$ws.Range("I122").Value = 555.9091
$ws.Range("J122").Value = 34066.332
$ws.Range("K122").Value = 5003.1819
$ws.Range("L122").Value = 306596.988
$ws.Range("M122").Value = -2553.1819
$ws.Range("N122").Value = -311496.988

$ws.Range("H131").Value = 1064.75
$ws.Range("J131").Value = 1064.75
$ws.Range("L131").Value = 3194.25
$ws.Range("N131").Value = -13274.25

$ws.Range("H132").Value = 3625.541
$ws.Range("J132").Value = 4327.0557
$ws.Range("L132").Value = 38943.5013
$ws.Range("N132").Value = -44003.5013

$ws.Range("H133").Value = 10779.167
$ws.Range("I133").Value = 4664.4443
$ws.Range("J133").Value = 16893.889
$ws.Range("K133").Value = 13993.3329
$ws.Range("L133").Value = 50681.667
$ws.Range("M133").Value = -8933.332900000001
$ws.Range("N133").Value = -60801.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 58.94737
$ws.Range("I2").Value = 62.1875
$ws.Range("J2").Value = 41.666668
$ws.Range("K2").Value = 62.1875
$ws.Range("L2").Value = 41.666668
$ws.Range("M2").Value = 50.8125
$ws.Range("N2").Value = -267.666668

$ws.Range("H11").Value = 415060.8
$ws.Range("I11").Value = 666766.7
$ws.Range("J11").Value = 37502
$ws.Range("K11").Value = 666766.7
$ws.Range("L11").Value = 37502
$ws.Range("M11").Value = -666627.7
$ws.Range("N11").Value = -37780

$ws.Range("H21").Value = 2883.3333
$ws.Range("J21").Value = 2883.3333
$ws.Range("L21").Value = 2883.3333
$ws.Range("N21").Value = -3229.3333

$ws.Range("H30").Value = 2883.3333
$ws.Range("J30").Value = 2883.3333
$ws.Range("L30").Value = 2883.3333
$ws.Range("N30").Value = -3093.3333

$ws.Range("H32").Value = 29999.334
$ws.Range("J32").Value = 29999.334
$ws.Range("L32").Value = 29999.334
$ws.Range("N32").Value = -30591.334

$ws.Range("H41").Value = 2609.4
$ws.Range("I41").Value = 2609.4
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2609.4
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2254.4
$ws.Range("N41").ClearContents()

$ws.Range("H51").Value = 47490
$ws.Range("J51").Value = 47490
$ws.Range("L51").Value = 47490
$ws.Range("N51").Value = -48508

$ws.Range("H80").Value = 674655.5600000001
$ws.Range("I80").Value = 1289711.2
$ws.Range("J80").Value = 59599.855
$ws.Range("K80").Value = 1289711.2
$ws.Range("L80").Value = 59599.855
$ws.Range("M80").Value = -1288713.2
$ws.Range("N80").Value = -61595.855

$ws.Range("H83").Value = 674655.5600000001
$ws.Range("I83").Value = 1289711.2
$ws.Range("J83").Value = 59599.855
$ws.Range("K83").Value = 6448556
$ws.Range("L83").Value = 297999.275
$ws.Range("M83").Value = -6443564
$ws.Range("N83").Value = -307983.275

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 40004.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 40004.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 40004.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -40344.5

$ws.Range("H31").Value = 1068.5714
$ws.Range("I31").Value = 1136.5
$ws.Range("J31").Value = 978
$ws.Range("K31").Value = 1136.5
$ws.Range("L31").Value = 978
$ws.Range("M31").Value = -888.5
$ws.Range("N31").Value = -1474

$ws.Range("H40").Value = 62503910
$ws.Range("I40").Value = 200001820
$ws.Range("J40").Value = 4864.091
$ws.Range("K40").Value = 200001820
$ws.Range("L40").Value = 4864.091
$ws.Range("M40").Value = -200001684
$ws.Range("N40").Value = -5136.091

$ws.Range("H82").Value = 2179.8
$ws.Range("I82").Value = 1966.6666
$ws.Range("J82").Value = 2499.5
$ws.Range("K82").Value = 1966.6666
$ws.Range("L82").Value = 2499.5
$ws.Range("M82").Value = -1605.6666
$ws.Range("N82").Value = -3221.5

$ws.Range("H85").Value = 2179.8
$ws.Range("I85").Value = 1966.6666
$ws.Range("J85").Value = 2499.5
$ws.Range("K85").Value = 1966.6666
$ws.Range("L85").Value = 2499.5
$ws.Range("M85").Value = -718.6666
$ws.Range("N85").Value = -4995.5

$ws.Range("H132").Value = 4617.722
$ws.Range("I132").Value = 4172.8
$ws.Range("J132").Value = 5173.875
$ws.Range("K132").Value = 12518.4
$ws.Range("L132").Value = 15521.625
$ws.Range("M132").Value = -9988.400000000001
$ws.Range("N132").Value = -20581.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 55504.5
$ws.Range("J10").Value = 55504.5
$ws.Range("L10").Value = 55504.5
$ws.Range("N10").Value = -55842.5

$ws.Range("H23").Value = 1200
$ws.Range("I23").Value = 1200
$ws.Range("K23").Value = 1200
$ws.Range("M23").Value = -971

$ws.Range("H86").Value = 122325
$ws.Range("J86").Value = 122325
$ws.Range("L86").Value = 122325
$ws.Range("N86").Value = -124571

$ws.Range("H89").Value = 122325
$ws.Range("J89").Value = 122325
$ws.Range("L89").Value = 611625
$ws.Range("N89").Value = -622857

$ws.Range("H107").Value = 1359.8
$ws.Range("I107").Value = 1399.5
$ws.Range("K107").Value = 4198.5
$ws.Range("M107").Value = -2278.5

$ws.Range("H122").Value = 2498.7693
$ws.Range("I122").Value = 2312.375
$ws.Range("K122").Value = 6937.125
$ws.Range("M122").Value = -4487.125
